# Auto-generated edit script applying the Lich_Profits market-data refresh.
# Updates currentAveragePrice* / Leve profit columns (H-N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# per the scheduled market-data runner, including a few rows whose price/profit
# cells collapsed to 0 (and had their now-meaningless profit cells cleared) and two
# rows that gained a new HQ-profit (column N) cell.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1075.6
$ws.Range("I19").Value = 596.2308
$ws.Range("K19").Value = 596.2308
$ws.Range("M19").Value = -421.2308
$ws.Range("H51").Value = 25713.715
$ws.Range("J51").Value = 25999.2
$ws.Range("L51").Value = 25999.2
$ws.Range("N51").Value = -26967.2
$ws.Range("H58").Value = 474.27274
$ws.Range("I58").Value = 502.125
$ws.Range("J58").Value = 400
$ws.Range("K58").Value = 1506.375
$ws.Range("L58").Value = 1200
$ws.Range("M58").Value = -1356.375
$ws.Range("N58").Value = -1500
$ws.Range("H92").Value = 2494.8975
$ws.Range("I92").Value = 2303.6667
$ws.Range("K92").Value = 2303.6667
$ws.Range("M92").Value = -1055.6667
$ws.Range("H99").Value = 258929000
$ws.Range("I99").Value = 17857224
$ws.Range("J99").Value = 500000740
$ws.Range("K99").Value = 53571672
$ws.Range("L99").Value = 1500002220
$ws.Range("M99").Value = -53570174
$ws.Range("N99").Value = -1500005216
$ws.Range("H113").Value = 5583.7437
$ws.Range("I113").Value = 6924.077
$ws.Range("K113").Value = 6924.077
$ws.Range("M113").Value = -3670.077
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 2625.3157
$ws.Range("I135").Value = 2127.6
$ws.Range("J135").Value = 3178.3333
$ws.Range("K135").Value = 19148.4
$ws.Range("L135").Value = 28604.9997
$ws.Range("M135").Value = -16613.4
$ws.Range("N135").Value = -33674.9997
$ws.Range("H137").Value = 3162874.5
$ws.Range("I137").Value = 5323493.5
$ws.Range("J137").Value = 5046.769
$ws.Range("K137").Value = 15970480.5
$ws.Range("L137").Value = 15140.307
$ws.Range("M137").Value = -15967930.5
$ws.Range("N137").Value = -20240.307
$ws.Range("N134").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1161.25
$ws.Range("I16").Value = 1148.3334
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1148.3334
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -861.3334
$ws.Range("N16").Value = -1774
$ws.Range("H32").Value = 7686.75
$ws.Range("I32").Value = 7153.56
$ws.Range("K32").Value = 7153.56
$ws.Range("M32").Value = -6866.56
$ws.Range("H61").Value = 4293.5747
$ws.Range("I61").Value = 2519
$ws.Range("K61").Value = 2519
$ws.Range("M61").Value = -2307
$ws.Range("H74").Value = 65036.16
$ws.Range("I74").Value = 67144.03
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 67144.03
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = -66270.03
$ws.Range("H77").Value = 65036.16
$ws.Range("I77").Value = 67144.03
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 335720.15
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = -331352.15
$ws.Range("H97").Value = 1853.1052
$ws.Range("I97").Value = 1743.3572
$ws.Range("K97").Value = 1743.3572
$ws.Range("M97").Value = -1247.3572
$ws.Range("H110").Value = 6708.9414
$ws.Range("I110").Value = 6876.9546
$ws.Range("J110").Value = 6400.9165
$ws.Range("K110").Value = 6876.9546
$ws.Range("L110").Value = 6400.9165
$ws.Range("M110").Value = -4831.9546
$ws.Range("N110").Value = -10490.9165
$ws.Range("H132").Value = 2871.75
$ws.Range("I132").Value = 2796.2856
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 8388.856800000001
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -5858.856800000001
$ws.Range("N132").Value = -15260
$ws.Range("H136").Value = 4293.5747
$ws.Range("I136").Value = 2519
$ws.Range("K136").Value = 7557
$ws.Range("M136").Value = -5007
$ws.Range("N74").Value = -3548
$ws.Range("N77").Value = -17736

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4198.4585
$ws.Range("I20").Value = 3384.1667
$ws.Range("J20").Value = 5012.75
$ws.Range("K20").Value = 3384.1667
$ws.Range("L20").Value = 5012.75
$ws.Range("M20").Value = -3137.1667
$ws.Range("N20").Value = -5506.75
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H134").Value = 2854
$ws.Range("J134").Value = 4995
$ws.Range("L134").Value = 14985
$ws.Range("N134").Value = -20055
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("H60").Value = 9103
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 9103
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 9103
$ws.Range("N60").Value = -10125
$ws.Range("H99").Value = 1264052.6
$ws.Range("I99").Value = 5000956
$ws.Range("K99").Value = 5000956
$ws.Range("M99").Value = -4999458
$ws.Range("H107").Value = 3588.7659
$ws.Range("I107").Value = 426.2857
$ws.Range("K107").Value = 426.2857
$ws.Range("M107").Value = 1493.7143
$ws.Range("H126").Value = 1264052.6
$ws.Range("I126").Value = 5000956
$ws.Range("K126").Value = 15002868
$ws.Range("M126").Value = -15000398
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("M60").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 115047.164
$ws.Range("I4").Value = 245345.19
$ws.Range("K4").Value = 736035.5700000001
$ws.Range("M4").Value = -735923.5700000001
$ws.Range("H18").Value = 547.5263
$ws.Range("I18").Value = 388
$ws.Range("K18").Value = 1164
$ws.Range("M18").Value = -995
$ws.Range("H32").Value = 28247.75
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("H39").Value = 1412.6333
$ws.Range("I39").Value = 899.05884
$ws.Range("K39").Value = 2697.17652
$ws.Range("M39").Value = -2403.17652
$ws.Range("H107").Value = 480.9565
$ws.Range("J107").Value = 483.22223
$ws.Range("L107").Value = 1449.66669
$ws.Range("N107").Value = -5289.66669
$ws.Range("H116").Value = 666
$ws.Range("I116").Value = 666
$ws.Range("K116").Value = 1998
$ws.Range("M116").Value = 1444
$ws.Range("H129").Value = 26517026
$ws.Range("J129").Value = 2360.4666
$ws.Range("L129").Value = 7081.399800000001
$ws.Range("N129").Value = -17081.3998
$ws.Range("M32").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 182.66667
$ws.Range("I2").Value = 199
$ws.Range("J2").Value = 174.5
$ws.Range("K2").Value = 199
$ws.Range("L2").Value = 174.5
$ws.Range("M2").Value = -86
$ws.Range("N2").Value = -400.5
$ws.Range("H80").Value = 4135.826
$ws.Range("I80").Value = 4272.8823
$ws.Range("J80").Value = 3747.5
$ws.Range("K80").Value = 4272.8823
$ws.Range("L80").Value = 3747.5
$ws.Range("M80").Value = -3274.8823
$ws.Range("N80").Value = -5743.5
$ws.Range("H83").Value = 4135.826
$ws.Range("I83").Value = 4272.8823
$ws.Range("J83").Value = 3747.5
$ws.Range("K83").Value = 21364.4115
$ws.Range("L83").Value = 18737.5
$ws.Range("M83").Value = -16372.4115
$ws.Range("N83").Value = -28721.5
$ws.Range("H107").Value = 652.2727
$ws.Range("I107").Value = 288.16666
$ws.Range("K107").Value = 288.16666
$ws.Range("M107").Value = 1631.83334
$ws.Range("H113").Value = 7791.5
$ws.Range("I113").Value = 3717.6667
$ws.Range("J113").Value = 20013
$ws.Range("K113").Value = 3717.6667
$ws.Range("L113").Value = 20013
$ws.Range("M113").Value = -1547.6667
$ws.Range("N113").Value = -24353
$ws.Range("H132").Value = 32637.428
$ws.Range("I132").Value = 33744.254
$ws.Range("K132").Value = 101232.762
$ws.Range("M132").Value = -98702.762
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 949999
$ws.Range("J139").Value = 949999
$ws.Range("L139").Value = 949999
$ws.Range("N139").Value = -960279
$ws.Range("N138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 50000576
$ws.Range("I16").Value = 55556140
$ws.Range("J16").Value = 510
$ws.Range("K16").Value = 55556140
$ws.Range("L16").Value = 510
$ws.Range("M16").Value = -55555970
$ws.Range("N16").Value = -850
$ws.Range("H46").Value = 1846.8
$ws.Range("I46").Value = 1996.2858
$ws.Range("J46").Value = 1498
$ws.Range("K46").Value = 1996.2858
$ws.Range("L46").Value = 1498
$ws.Range("M46").Value = -1808.2858
$ws.Range("N46").Value = -1874
$ws.Range("H122").Value = 8315.040000000001
$ws.Range("I122").Value = 8019.2
$ws.Range("K122").Value = 24057.6
$ws.Range("M122").Value = -21607.6
$ws.Range("H136").Value = 4403.6665
$ws.Range("I136").Value = 2274.4
$ws.Range("J136").Value = 15050
$ws.Range("K136").Value = 6823.200000000001
$ws.Range("L136").Value = 45150
$ws.Range("M136").Value = -4273.200000000001
$ws.Range("N136").Value = -50250

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 298626.03
$ws.Range("I136").Value = 421653.75
$ws.Range("J136").Value = 3359.5
$ws.Range("K136").Value = 1264961.25
$ws.Range("L136").Value = 10078.5
$ws.Range("M136").Value = -1262411.25
$ws.Range("N136").Value = -15178.5
$ws.Range("H139").Value = 89786
$ws.Range("J139").Value = 89786
$ws.Range("L139").Value = 89786
$ws.Range("N139").Value = -100066
